$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.554.90'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.84%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.745.65'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.69%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '116.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.33%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '331.56'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.24%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.533'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.40%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.05%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.566'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.02%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.78'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.76%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.40'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.97%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0832'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.35%  '

# Row 13
$ws.Range('E13').Value = '  +2.61%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.69'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.63%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.168.13'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.57%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.721.91'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.17%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.886'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.60%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.430.90'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.76%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.74'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.99%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.05'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.75%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.86'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.41%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0964'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.06%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '287.32'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.94%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.65'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.07%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.61'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.29%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.96'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.23%  '

# Row 27
$ws.Range('E27').Value = '  +0.02%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.33'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.44%  '

# Row 29
$ws.Range('E29').Value = '  -0.64%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.142'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.46%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.81'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.95%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '50.08'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.70%  '

# Row 33
$ws.Range('E33').Value = '  +1.26%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0830'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.36%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '19.45'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.88%  '

# Row 36
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.06'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.17%  '

# Row 37
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.05%  '

# Row 38
$ws.Range('E38').Value = '  +1.24%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.22'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.79%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.78'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.92%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '129.23'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.34%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0354'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +10.46%  '

# Row 43
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.29'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.99%  '

# Row 44
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.113'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.33%  '

# Row 45
$ws.Range('E45').Value = '  +1.88%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.115.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.03%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.25'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +12.13%  '

# Row 48
$ws.Range('E48').Value = '  -2.13%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.54'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.66%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.12'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.27%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '60.36'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.36%  '
